# "Fixed foil for this.PI"
#
# Slide 10 ("Writing Your Own Module" - module.exports code sample):
#   - "...return PI * r * r;"    -> "...return this.PI * r * r;"
#   - "...return 2 * PI * r;"    -> "...return 2 * this.PI * r;"
#
# Slide 12 / Slide 3: the underlying wording doesn't change, but the
# runs that made up a couple of sentences get re-typed/merged into a
# single run (as happens when someone retypes a sentence in the UI
# instead of leaving the old run boundaries in place). We reproduce
# that by rewriting those paragraphs/sub-ranges through TextRange so
# the runtime re-flows them into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10: "return PI" -> "return this.PI" and "2 * PI" -> "2 * this.PI"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$code = $s10.Shapes.Item(1).TextFrame.TextRange

# Paragraph 4: "    return PI * r * r;" -- replace the "PI" in "return PI "
$para4 = $code.Paragraphs(4, 1)
$para4.Characters(12, 2).Text = "this.PI"

# Paragraph 7: "    return 2 * PI * r;" -- replace the "PI" in "2 * PI "
$para7 = $code.Paragraphs(7, 1)
$para7.Characters(16, 2).Text = "this.PI"

# ---------------------------------------------------------------------------
# Slide 12 ("Summary"): collapse runs that were split mid-sentence back
# into single runs (text content is unchanged).
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$sum = $s12.Shapes.Item(1).TextFrame.TextRange

# Paragraph 1: "Modules: " + "self" + "-contained behavior"
# (Use a throwaway placeholder with no char in common with the real text
# so the runtime's text-diff can't fold the edit back into multiple runs.)
$sum.Paragraphs(1, 1).Text = "@"
$sum.Paragraphs(1, 1).Text = "Modules: self-contained behavior"

# Paragraph 2: "...(" + "npm" + "): " + "installs " + "& uninstalls modules"
# Only merge the tail ("): installs & uninstalls modules"); leave
# "Node Package Manager (" and "npm" alone.
$p2 = $sum.Paragraphs(2, 1)
$p2.Characters(26, 32).Text = "@"
$sum.Paragraphs(2, 1).Characters(26, 1).Text = "): installs & uninstalls modules"

# Paragraph 3: "Using modules: " + "use " + "require('module')" (Monaco)
# Merge the first two runs into "Using modules: use "; leave the Monaco
# run alone.
$p3 = $sum.Paragraphs(3, 1)
$p3.Characters(1, 19).Text = "@"
$sum.Paragraphs(3, 1).Characters(1, 1).Text = "Using modules: use "

# Paragraph 4: "Writing modules" + ": " + "remember " + "index.js" + ...
# Merge the first three runs into "Writing modules: remember "; leave
# the rest alone.
$p4 = $sum.Paragraphs(4, 1)
$p4.Characters(1, 26).Text = "@"
$sum.Paragraphs(4, 1).Characters(1, 1).Text = "Writing modules: remember "

# ---------------------------------------------------------------------------
# Slide 3 ("What is a Module?"): collapse "To " + "export an artifact..."
# into a single run (text content is unchanged).
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$what = $s3.Shapes.Item(1).TextFrame.TextRange

$p5 = $what.Paragraphs(5, 1)
$p5.Characters(1, 48).Text = "@"
$what.Paragraphs(5, 1).Characters(1, 1).Text = "To export an artifact from a module, use Node's "
